$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dynamic Programming section: mark House Robber II (row 28), Decode Ways
# (row 29), Unique Paths (row 30) and Jump Game (row 31) as completed by
# flipping column C from the "<yes/no>" placeholder to "<yes>".
$ws.Range("C28").Value = "<yes>"
$ws.Range("C29").Value = "<yes>"
$ws.Range("C30").Value = "<yes>"
$ws.Range("C31").Value = "<yes>"

# Leave the view scrolled/selected where the author ended up after making
# the edits: viewport anchored at row 19, cursor resting on C31.
$excel.Goto($ws.Range("A19"), $true)
$ws.Range("C31").Select()
